# Update cryptos list (Price / Volume(1h) columns, and two coin-name swaps)
# as published on Sat May 13 03:50:02 UTC 2023 with GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value2 = '26.910.34'
$ws.Range('E2').Value2 = '  +0.20%  '
# Row 3 - Ethereum
$ws.Range('D3').Value2 = '1.812.70'
$ws.Range('E3').Value2 = '  +1.62%  '
# Row 4 - TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.002'
$ws.Range('E4').Value2 = '  -0.64%  '
# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '309.45'
$ws.Range('E5').Value2 = '  +0.45%  '
# Row 6 - USDC
$ws.Range('E6').Value2 = '  -0.54%  '
# Row 7 - XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.4286'
$ws.Range('E7').Value2 = '  +1.46%  '
# Row 8 - Cardano
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.3688'
# Row 9 - Dogecoin
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.07232'
$ws.Range('E9').Value2 = '  +0.97%  '
# Row 10 - Polygon
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.8613'
$ws.Range('E10').Value2 = '  +2.76%  '
# Row 11 - Solana
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '21.18'
$ws.Range('E11').Value2 = '  +4.66%  '
# Row 12 - WrappedEther
$ws.Range('D12').Value2 = '2.012.18'
$ws.Range('E12').Value2 = '  +9.53%  '
# Row 13 - Chainlink
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '6.632'
$ws.Range('E13').Value2 = '  +4.76%  '
# Row 14 - Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '5.392'
$ws.Range('E14').Value2 = '  +2.84%  '
# Row 15 - TRON
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.06888'
$ws.Range('E15').Value2 = '  +1.23%  '
# Row 16 - Litecoin
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '80.61'
$ws.Range('E16').Value2 = '  +1.88%  '
# Row 17 - BinanceUSD
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '1.003'
$ws.Range('E17').Value2 = '  -0.84%  '
# Row 18 - ShibaInu
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '0.000008861'
$ws.Range('E18').Value2 = '  +2.07%  '
# Row 19 - Dai
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '1.003'
$ws.Range('E19').Value2 = '  -0.51%  '
# Row 20 - Avalanche
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '15.18'
$ws.Range('E20').Value2 = '  +1.84%  '
# Row 21 - WrappedBTC
$ws.Range('D21').Value2 = '26.949.06'
$ws.Range('E21').Value2 = '  -0.34%  '
# Row 22 - Uniswap
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '5.176'
$ws.Range('E22').Value2 = '  +2.89%  '
# Row 23 - Cosmos
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '11.03'
$ws.Range('E23').Value2 = '  +0.07%  '
# Row 24 - WrappedliquidstakedEther2.0
$ws.Range('D24').Value2 = '2.243.84'
$ws.Range('E24').Value2 = '  +9.76%  '
# Row 25 - Monero
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '153.60'
$ws.Range('E25').Value2 = '  +0.46%  '
# Row 26 - Toncoin
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '1.881'
$ws.Range('E26').Value2 = '  -1.95%  '
# Row 27 - EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '18.25'
$ws.Range('E27').Value2 = '  +0.82%  '
# Row 28 - InternetComputer(DFINITY)
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '5.230'
$ws.Range('E28').Value2 = '  +4.22%  '
# Row 29 - BitcoinCash
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '114.79'
$ws.Range('E29').Value2 = '  +0.58%  '
# Row 30 - LidoDAOToken
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '1.877'
$ws.Range('E30').Value2 = '  +15.13%  '
# Row 31 - Stellar
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '0.08947'
$ws.Range('E31').Value2 = '  +0.13%  '
# Row 32 - ImmutableX -> ARBITRUM (swap)
$ws.Range('B32').Value2 = 'ARBITRUM'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '1.166'
$ws.Range('E32').Value2 = '  +7.41%  '
# Row 33 - ARBITRUM -> ImmutableX (swap)
$ws.Range('B33').Value2 = 'ImmutableX'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '0.7432'
$ws.Range('E33').Value2 = '  +4.34%  '
# Row 34 - Filecoin
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '4.421'
$ws.Range('E34').Value2 = '  +2.71%  '
# Row 35 - HuobiToken
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '2.794'
$ws.Range('E35').Value2 = '  -1.74%  '
# Row 36 - Frax
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '1.006'
$ws.Range('E36').Value2 = '  -0.18%  '
# Row 37 - TrustWalletToken
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '1.122'
$ws.Range('E37').Value2 = '  +4.22%  '
# Row 38 - Hedera
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.05209'
$ws.Range('E38').Value2 = '  +2.52%  '
# Row 39 - VeChain
$ws.Range('E39').Value2 = '  +1.49%  '
# Row 40 - TheSandbox
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.5086'
$ws.Range('E40').Value2 = '  +3.27%  '
# Row 41 - Algorand
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.1641'
$ws.Range('E41').Value2 = '  +2.00%  '
# Row 42 - MXToken
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '2.710'
$ws.Range('E42').Value2 = '  +7.77%  '
# Row 43 - FraxShare
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '6.431'
$ws.Range('E43').Value2 = '  +7.40%  '
# Row 44 - Aptos
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '8.283'
$ws.Range('E44').Value2 = '  +4.80%  '
# Row 45 - Quant
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '106.74'
$ws.Range('E45').Value2 = '  +2.49%  '
# Row 46 - EnergySwap
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '10.40'
$ws.Range('E46').Value2 = '  +3.80%  '
# Row 47 - PaxDollar
$ws.Range('E47').Value2 = '  -0.44%  '
# Row 48 - NEARProtocol
$ws.Range('E48').Value2 = '  +5.10%  '
# Row 49 - Cronos -> Decentraland (swap)
$ws.Range('B49').Value2 = 'Decentraland'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '0.4570'
$ws.Range('E49').Value2 = '  +2.38%  '
# Row 50 - Decentraland -> Cronos (swap)
$ws.Range('B50').Value2 = 'Cronos'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.06274'
$ws.Range('E50').Value2 = '  +0.40%  '
# Row 51 - RenderToken
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '1.807'
$ws.Range('E51').Value2 = '  +5.98%  '
